$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Bánh bi" -> "Bi mắt trâu" for item 10 (row 11) ---
$ws.Range("B11").Value = "Bi mắt trâu"

# --- Row 11: quantity / price / total changed (2 x 20000 = 40000) ---
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 20000

# --- Remove the old SUM(E:E) total that lived in G2, keep its style ---
$ws.Range("G2").ClearContents()

# --- Add a "Tổng" / grand-total box in I6:I7 ---
$ws.Range("I6:I7").Font.Color = 255
$ws.Range("I6:I7").Interior.Color = 65535
$ws.Range("I7").Borders.Weight = -4138

$ws.Range("I6").Value = "Tổng"
$ws.Range("I7").Formula = "=SUM(E:E)"

# --- Relabel the E1 header from "Tổng" to "Thành tiền" ---
$ws.Range("E1").Value = "Thành tiền"

# --- Selection moved ---
$ws.Range("F19").Select()
